$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "2024-07-03"
$ws.Range("B2").Value = "하스"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 289.6
$ws.Range("E2").Value = "삼성"
$ws.Range("F2").Value = 289.6
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 16000
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "2024-06-24"
$ws.Range("P2").Value = "2024-06-27"
$ws.Range("Q2").Value = 1357500
